$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price/Volume columns so numeric-looking
# strings (e.g. "357.25", "35.00", "0.999") are preserved verbatim as text,
# matching the source workbook which stores these as inline strings.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.763.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.778.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.25"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.31"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.557"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.587"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.82"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0846"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.52"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.62"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.210.48"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.785.08"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.935"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.679.12"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.08"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.18"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0971"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.22"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.71"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.41"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +16.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.25"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.27"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.73%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.80"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0453"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -8.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0840"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.78"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.14"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.04%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.23"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.41"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.74"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.087.45"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "SEI"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.947"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.61%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.65"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.04%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.60%  "
